# Auto-generated: apply scheduled-runner market data update to Sheets/Sagittarius_Profits.xlsx
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) across multiple job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 478.57144
$ws.Range("I8").Value = 92.5
$ws.Range("K8").Value = 277.5
$ws.Range("M8").Value = -138.5
$ws.Range("H53").Value = 315.66666
$ws.Range("I53").Value = 216.1
$ws.Range("J53").Value = 514.8
$ws.Range("K53").Value = 216.1
$ws.Range("L53").Value = 514.8
$ws.Range("M53").Value = 420.9
$ws.Range("N53").Value = -1788.8
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("H98").Value = 1912
$ws.Range("I98").Value = 1912
$ws.Range("K98").Value = 1912
$ws.Range("M98").Value = -414
$ws.Range("H122").Value = 1912
$ws.Range("I122").Value = 1912
$ws.Range("K122").Value = 5736
$ws.Range("M122").Value = -3286
$ws.Range("H132").Value = 2919.7
$ws.Range("I132").Value = 2937.125
$ws.Range("K132").Value = 8811.375
$ws.Range("M132").Value = -6281.375
$ws.Range("H135").Value = 1715.6666
$ws.Range("I135").Value = 1793.4667
$ws.Range("K135").Value = 16141.2003
$ws.Range("M135").Value = -13606.2003
$ws.Range("H137").Value = 1132
$ws.Range("I137").Value = 903.7692
$ws.Range("J137").Value = 1555.8572
$ws.Range("K137").Value = 2711.3076
$ws.Range("L137").Value = 4667.571599999999
$ws.Range("M137").Value = -161.3076000000001
$ws.Range("N137").Value = -9767.571599999999
$ws.Range("H138").Value = 5526.1113
$ws.Range("I138").Value = 5057.2
$ws.Range("K138").Value = 15171.6
$ws.Range("M138").Value = -10031.6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 15079
$ws.Range("I36").Value = 5118.5
$ws.Range("K36").Value = 5118.5
$ws.Range("M36").Value = -4772.5
$ws.Range("H45").Value = 4597.9
$ws.Range("I45").Value = 4597.9
$ws.Range("K45").Value = 4597.9
$ws.Range("M45").Value = -4220.9
$ws.Range("H61").Value = 1787.4445
$ws.Range("I61").Value = 1819.8667
$ws.Range("K61").Value = 1819.8667
$ws.Range("M61").Value = -1607.8667
$ws.Range("H74").Value = 1354.6666
$ws.Range("I74").Value = 1354.6666
$ws.Range("K74").Value = 1354.6666
$ws.Range("M74").Value = -480.6666
$ws.Range("H77").Value = 1354.6666
$ws.Range("I77").Value = 1354.6666
$ws.Range("K77").Value = 6773.333000000001
$ws.Range("M77").Value = -2405.333000000001
$ws.Range("H132").Value = 1233
$ws.Range("I132").Value = 1233
$ws.Range("K132").Value = 3699
$ws.Range("M132").Value = -1169
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 1787.4445
$ws.Range("I136").Value = 1819.8667
$ws.Range("K136").Value = 5459.6001
$ws.Range("M136").Value = -2909.6001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2577.5
$ws.Range("I105").Value = 2577.5
$ws.Range("K105").Value = 2577.5
$ws.Range("M105").Value = -830.5
$ws.Range("H107").Value = 956.8
$ws.Range("I107").Value = 1011.3333
$ws.Range("J107").Value = 875
$ws.Range("K107").Value = 1011.3333
$ws.Range("L107").Value = 875
$ws.Range("M107").Value = 908.6667
$ws.Range("N107").Value = -4715
$ws.Range("H134").Value = 2588.3333
$ws.Range("I134").Value = 2686.875
$ws.Range("K134").Value = 8060.625
$ws.Range("M134").Value = -5525.625
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 507
$ws.Range("I10").Value = 507
$ws.Range("K10").Value = 507
$ws.Range("M10").Value = -368
$ws.Range("H31").Value = 2317.2593
$ws.Range("I31").Value = 1803.6666
$ws.Range("J31").Value = 2574.0557
$ws.Range("K31").Value = 1803.6666
$ws.Range("L31").Value = 2574.0557
$ws.Range("M31").Value = -1508.6666
$ws.Range("N31").Value = -3164.0557
$ws.Range("H34").Value = 2317.2593
$ws.Range("I34").Value = 1803.6666
$ws.Range("J34").Value = 2574.0557
$ws.Range("K34").Value = 1803.6666
$ws.Range("L34").Value = 2574.0557
$ws.Range("M34").Value = -1601.6666
$ws.Range("N34").Value = -2978.0557
$ws.Range("H58").Value = 4388.25
$ws.Range("I58").Value = 3258.6667
$ws.Range("K58").Value = 3258.6667
$ws.Range("M58").Value = -3055.6667
$ws.Range("H59").Value = 130250
$ws.Range("I59").Value = 1000
$ws.Range("J59").Value = 173333.33
$ws.Range("K59").Value = 1000
$ws.Range("L59").Value = 173333.33
$ws.Range("M59").Value = 145
$ws.Range("N59").Value = -175623.33
$ws.Range("H132").Value = 1509.125
$ws.Range("I132").Value = 1509.125
$ws.Range("K132").Value = 4527.375
$ws.Range("M132").Value = -1997.375
$ws.Range("H136").Value = 4388.25
$ws.Range("I136").Value = 3258.6667
$ws.Range("K136").Value = 9776.000100000001
$ws.Range("M136").Value = -7226.000100000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 159.72728
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = 172.7
$ws.Range("K2").Value = 180
$ws.Range("L2").Value = 1036.2
$ws.Range("M2").Value = -67
$ws.Range("N2").Value = -1262.2
$ws.Range("H4").Value = 14993074
$ws.Range("I4").Value = 18336796
$ws.Range("J4").Value = 4961905
$ws.Range("K4").Value = 55010388
$ws.Range("L4").Value = 14885715
$ws.Range("M4").Value = -55010276
$ws.Range("N4").Value = -14885939
$ws.Range("H121").Value = 3327.625
$ws.Range("I121").Value = 675
$ws.Range("K121").Value = 2025
$ws.Range("M121").Value = -715
$ws.Range("H122").Value = 393
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 393
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 3537
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -8437
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 79199.8
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 79199.8
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 79199.8
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -80839.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3360.3
$ws.Range("I46").Value = 1826.6
$ws.Range("K46").Value = 1826.6
$ws.Range("M46").Value = -1638.6
$ws.Range("H68").Value = 1812.375
$ws.Range("I68").Value = 1399.6
$ws.Range("J68").Value = 2500.3333
$ws.Range("K68").Value = 1399.6
$ws.Range("L68").Value = 2500.3333
$ws.Range("M68").Value = -650.5999999999999
$ws.Range("N68").Value = -3998.3333
$ws.Range("H71").Value = 1812.375
$ws.Range("I71").Value = 1399.6
$ws.Range("J71").Value = 2500.3333
$ws.Range("K71").Value = 6998
$ws.Range("L71").Value = 12501.6665
$ws.Range("M71").Value = -3254
$ws.Range("N71").Value = -19989.6665
$ws.Range("H136").Value = 4040
$ws.Range("I136").Value = 3567
$ws.Range("J136").Value = 4749.5
$ws.Range("K136").Value = 10701
$ws.Range("L136").Value = 14248.5
$ws.Range("M136").Value = -8151
$ws.Range("N136").Value = -19348.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7609.1665
$ws.Range("J122").Value = 7917.75
$ws.Range("L122").Value = 23753.25
$ws.Range("N122").Value = -28653.25
$ws.Range("H126").Value = 1320.3846
$ws.Range("I126").Value = 1386.6
$ws.Range("J126").Value = 1099.6666
$ws.Range("K126").Value = 4159.799999999999
$ws.Range("L126").Value = 3298.9998
$ws.Range("M126").Value = -1689.799999999999
$ws.Range("N126").Value = -8238.9998
$ws.Range("H132").Value = 533
$ws.Range("I132").Value = 533
$ws.Range("K132").Value = 1599
$ws.Range("M132").Value = 931
$ws.Range("H136").Value = 2352.0527
$ws.Range("I136").Value = 2215.6924
$ws.Range("K136").Value = 6647.0772
$ws.Range("M136").Value = -4097.0772
